$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the Goods table (rows 2-21) ---
# Row 2 (food1)
$ws.Range("D2").Value = 80

# Row 3 (food2)
$ws.Range("B3").Value = 20
$ws.Range("D3").Value = 12

# Row 4 (food3)
$ws.Range("B4").Value = 86
$ws.Range("D4").Value = 70

# Row 5 (food4)
$ws.Range("B5").Value = 107
$ws.Range("D5").Value = 80

# Row 6 (food5)
$ws.Range("B6").Value = 5
$ws.Range("D6").Value = 3

# Row 7 (was food6 -> item1)
$ws.Range("A7").Value = "item1"
$ws.Range("B7").Value = 34
$ws.Range("D7").Value = 20

# Row 8 (was food7 -> item2)
$ws.Range("A8").Value = "item2"
$ws.Range("D8").Value = 90

# Row 9 (was food8 -> item3)
$ws.Range("A9").Value = "item3"
$ws.Range("B9").Value = 255
$ws.Range("D9").Value = 180

# Row 10 (was food9 -> item4)
$ws.Range("A10").Value = "item4"
$ws.Range("B10").Value = 12
$ws.Range("D10").Value = 8

# Row 11 (was food10 -> item5)
$ws.Range("A11").Value = "item5"
$ws.Range("B11").Value = 48
$ws.Range("D11").Value = 35

# Row 12 (was food11 -> necessities1)
$ws.Range("A12").Value = "necessities1"
$ws.Range("B12").Value = 8
$ws.Range("D12").Value = 5

# Row 13 (was food12 -> necessities2)
$ws.Range("A13").Value = "necessities2"
$ws.Range("B13").Value = 26
$ws.Range("D13").Value = 20

# Row 14 (was food13 -> necessities3)
$ws.Range("A14").Value = "necessities3"
$ws.Range("B14").Value = 74
$ws.Range("D14").Value = 55

# Row 15 (was food14 -> necessities4)
$ws.Range("A15").Value = "necessities4"
$ws.Range("B15").Value = 39
$ws.Range("D15").Value = 30

# Row 16 (was food15 -> necessities5)
$ws.Range("A16").Value = "necessities5"
$ws.Range("B16").Value = 101
$ws.Range("D16").Value = 80

# Row 17 (was food16 -> tools1)
$ws.Range("A17").Value = "tools1"
$ws.Range("B17").Value = 350
$ws.Range("D17").Value = 280

# Row 18 (was food17 -> tools2)
$ws.Range("A18").Value = "tools2"
$ws.Range("B18").Value = 110
$ws.Range("D18").Value = 90

# Row 19 (was food18 -> tools3)
$ws.Range("A19").Value = "tools3"
$ws.Range("B19").Value = 80
$ws.Range("D19").Value = 50

# Row 20 (was food19 -> tools4)
$ws.Range("A20").Value = "tools4"
$ws.Range("B20").Value = 500
$ws.Range("D20").Value = 280

# Row 21 (was food20 -> tools5)
$ws.Range("A21").Value = "tools5"
$ws.Range("B21").Value = 30
$ws.Range("D21").Value = 15

# --- Update the view: selection moves from F11 to D21 ---
$ws.Range("D21").Select()
